$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "Area" / "Atotal" columns (G, H) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# --- B2 / C2 used to hold the placeholder string "-"; they become numeric 0 ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New "Area" column G: incremental cross-sectional area per segment ---
# G2 is special-cased against an implicit 0 baseline (no row above it)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
# G3 also stands on its own (mirrors the pattern of the existing Q/E column)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
# G4:G15 share one formula, filled down
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- New "Atotal" column H: running total of the Area column ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Row 14 previously ended the sheet at column A only; row 15 is brand new ---
$ws.Range("A14").Value = 210

# --- Selection moves to the newly added Area/Atotal block ---
[void]$ws.Range("G1:H15").Select()
